$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "cryptos" price list (Coin / Link / Price / Volume(1h)) with the
# latest scrape. Most rows only get refreshed Price / Volume(1h) figures;
# rows 35-37 also got reshuffled (Binance-Peg BSC-USD, Aptos and
# RenzoRestakedETH swapped ranking order). Price cells whose new text would
# otherwise be auto-parsed as a plain decimal number are first switched to
# Text format ("@") so Excel keeps storing them as the literal price string
# (e.g. trailing zeros such as "1.00" or "163.60" survive) instead of
# silently collapsing them into a number.
$ws.Range("D2").Value = "70.910.83"
$ws.Range("E2").Value = "  -0.98%  "
$ws.Range("D3").Value = "3.826.83"
$ws.Range("E3").Value = "  +0.49%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "702.17"
$ws.Range("E5").Value = "  -0.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.92"
$ws.Range("E6").Value = "  -1.44%  "
$ws.Range("D7").Value = "3.825.50"
$ws.Range("E7").Value = "  +0.51%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -0.73%  "
$ws.Range("E10").Value = "  -0.64%  "
$ws.Range("E11").Value = "  -1.31%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.458"
$ws.Range("E12").Value = "  -0.63%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000254"
$ws.Range("E13").Value = "  -1.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.68"
$ws.Range("E14").Value = "  -0.13%  "
$ws.Range("D15").Value = "4.474.26"
$ws.Range("E15").Value = "  +0.44%  "
$ws.Range("D16").Value = "3.890.85"
$ws.Range("E16").Value = "  +2.28%  "
$ws.Range("D17").Value = "70.937.84"
$ws.Range("E17").Value = "  -0.90%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.21"
$ws.Range("E18").Value = "  -0.46%  "
$ws.Range("E19").Value = "  +0.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.37"
$ws.Range("E20").Value = "  -2.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "496.18"
$ws.Range("E21").Value = "  +2.42%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.65"
$ws.Range("E22").Value = "  -3.99%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.732"
$ws.Range("E23").Value = "  +2.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.17"
$ws.Range("E24").Value = "  +0.59%  "
$ws.Range("E25").Value = "  +0.20%  "
$ws.Range("E26").Value = "  +0.31%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.13"
$ws.Range("E27").Value = "  -1.75%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.08"
$ws.Range("E28").Value = "  -3.92%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.06"
$ws.Range("E30").Value = "  -2.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.44"
$ws.Range("E31").Value = "  -1.88%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.22"
$ws.Range("E32").Value = "  -4.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.26"
$ws.Range("E33").Value = "  -1.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.177"
$ws.Range("E34").Value = "  -4.28%  "
$ws.Range("B35").Value = "Binance-PegBSC-USD"
$ws.Range("C35").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.01"
$ws.Range("E35").Value = "  +1.03%  "
$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.18"
$ws.Range("E36").Value = "  -1.50%  "
$ws.Range("B37").Value = "RenzoRestakedETH"
$ws.Range("C37").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D37").Value = "3.792.62"
$ws.Range("E37").Value = "  +0.81%  "
$ws.Range("E38").Value = "  -1.25%  "
$ws.Range("E39").Value = "  -2.94%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.03"
$ws.Range("E40").Value = "  +4.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.96"
$ws.Range("E41").Value = "  -1.48%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.32"
$ws.Range("E42").Value = "  -3.94%  "
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.000315"
$ws.Range("E45").Value = "  +1.37%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "163.60"
$ws.Range("E46").Value = "  -0.40%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "431.53"
$ws.Range("E47").Value = "  +3.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "48.91"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.70"
$ws.Range("E49").Value = "  +0.26%  "
$ws.Range("E50").Value = "  -0.30%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.296"
$ws.Range("E51").Value = "  -2.46%  "
